$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20 ("Invert Binary Tree") re-uses the exact same per-row style as every
# other "Easy" row (e.g. row 19) - copy its formatting down first, then fill
# in the values, so the new cells land on the existing shared xfs instead of
# inheriting the bare column defaults.
$ws.Range("A19:G19").Copy()
$ws.Range("A20:G20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A20").Value = 226
$ws.Range("B20").Value = "Easy"
$ws.Range("C20").Value = "Invert Binary Tree"
$ws.Range("D20").Value = "https://shorturl.at/Lq38f"
$ws.Range("E20").Value = "Recursion"
$ws.Range("F20").Value = "O(n)"
$ws.Range("G20").Value = "Simple recursion"

$ws.Hyperlinks.Add($ws.Range("D20"), "https://shorturl.at/Lq38f")

# Adding the hyperlink re-stamps D20's cell format; reapply D19's formatting
# so D20 collapses back onto the shared "hyperlink" style used by the rest
# of the url column instead of keeping a one-off duplicate.
$ws.Range("D19").Copy()
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("C22").Select()
